$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CT 01): test re-run as "No" with no status recorded, new output timestamp
$ws.Range("B2").Value = "No"
$ws.Range("C2").Value = "'"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Value = "25/05/2020"

# Row 3 (CT 02): test re-run as "Yes", new output timestamp
$ws.Range("B3").Value = "Yes"
$ws.Range("H3").Value = "25/05/2020"

# Row 4 (CT 03): test re-run as "No" with no status recorded, new output timestamp
$ws.Range("B4").Value = "No"
$ws.Range("C4").Value = "'"
$ws.Range("D4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Value = "25/05/2020"

# Row 6 (CT 05): test re-run as "No" with no status recorded, new output timestamp
$ws.Range("B6").Value = "No"
$ws.Range("C6").Value = "'"
$ws.Range("D6").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("H6").Value = "25/05/2020"

# Update selection to reflect the new active cell
$ws.Range("B10").Select()
